$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the greeting text in E8 from "Good Morning" to "GIT UPDATE"
$ws.Range("E8").Value = "GIT UPDATE"

# Make E8 the active/selected cell (as captured in the sheet view)
$ws.Range("E8").Select()
